$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.Style = "Normal"
}

# --- Update Price (D) / Volume(1h) (E) values for existing rows ---

Set-TextValue $ws.Range('D2') '26.738.65'
Set-TextValue $ws.Range('E2') '  +0.30%  '

Set-TextValue $ws.Range('D3') '1.601.94'
Set-TextValue $ws.Range('E3') '  +0.28%  '

Set-TextValue $ws.Range('E4') '  +0.17%  '

Set-TextValue $ws.Range('D5') '211.73'
Set-TextValue $ws.Range('E5') '  +0.13%  '

Set-TextValue $ws.Range('E6') '  -0.10%  '

Set-TextValue $ws.Range('D8') '0.0619'
Set-TextValue $ws.Range('E8') '  +0.06%  '

Set-TextValue $ws.Range('E9') '  -0.27%  '

Set-TextValue $ws.Range('D10') '19.68'
Set-TextValue $ws.Range('E10') '  +0.57%  '

Set-TextValue $ws.Range('E11') '  +0.73%  '

Set-TextValue $ws.Range('D12') '1.826.00'
Set-TextValue $ws.Range('E12') '  +0.24%  '

Set-TextValue $ws.Range('D13') '1.602.28'
Set-TextValue $ws.Range('E13') '  -0.09%  '

Set-TextValue $ws.Range('E14') '  +0.61%  '

Set-TextValue $ws.Range('D15') '0.524'
Set-TextValue $ws.Range('E15') '  +0.01%  '

Set-TextValue $ws.Range('D16') '65.05'
Set-TextValue $ws.Range('E16') '  -0.02%  '

Set-TextValue $ws.Range('D18') '209.84'
Set-TextValue $ws.Range('E18') '  +0.37%  '

Set-TextValue $ws.Range('E19') '  +0.19%  '

Set-TextValue $ws.Range('D20') '7.17'
Set-TextValue $ws.Range('E20') '  +2.12%  '

Set-TextValue $ws.Range('E21') '  +0.16%  '

Set-TextValue $ws.Range('D22') '2.23'
Set-TextValue $ws.Range('E22') '  -3.54%  '

Set-TextValue $ws.Range('E23') '  -0.09%  '

Set-TextValue $ws.Range('D24') '143.50'
Set-TextValue $ws.Range('E24') '  -0.45%  '

Set-TextValue $ws.Range('E25') '  +0.38%  '

Set-TextValue $ws.Range('D26') '7.08'
Set-TextValue $ws.Range('E26') '  -0.64%  '

Set-TextValue $ws.Range('E27') '  -0.63%  '

Set-TextValue $ws.Range('D28') '15.33'
Set-TextValue $ws.Range('E28') '  +0.29%  '

Set-TextValue $ws.Range('D29') '0.0509'
Set-TextValue $ws.Range('E29') '  -1.14%  '

Set-TextValue $ws.Range('E30') '  +0.31%  '

Set-TextValue $ws.Range('E31') '  +0.52%  '

Set-TextValue $ws.Range('E32') '  +0.28%  '

Set-TextValue $ws.Range('D33') '1.291.57'
Set-TextValue $ws.Range('E33') '  +0.21%  '

Set-TextValue $ws.Range('E34') '  +0.57%  '

Set-TextValue $ws.Range('E35') '  +0.33%  '

Set-TextValue $ws.Range('E36') '  -2.18%  '

Set-TextValue $ws.Range('E37') '  +11.34%  '

Set-TextValue $ws.Range('E38') '  -0.05%  '

Set-TextValue $ws.Range('E39') '  -0.11%  '

Set-TextValue $ws.Range('E40') '  -1.50%  '

Set-TextValue $ws.Range('E41') '  -0.68%  '

Set-TextValue $ws.Range('D42') '0.783'
Set-TextValue $ws.Range('E42') '  -0.10%  '

Set-TextValue $ws.Range('D43') '62.98'
Set-TextValue $ws.Range('E43') '  -0.88%  '

Set-TextValue $ws.Range('D44') '1.738.00'
Set-TextValue $ws.Range('E44') '  +0.22%  '

Set-TextValue $ws.Range('D45') '90.47'
Set-TextValue $ws.Range('E45') '  -0.27%  '

# --- Rows 47-51: BabyDogeCoin (old row 47) is removed from the listing,
# the remaining coins shift up one row, and a new coin (Mantle) is
# appended as the new last row (51). ---

Set-TextValue $ws.Range('B47') 'Algorand'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D47') '0.102'
Set-TextValue $ws.Range('E47') '  +0.22%  '

Set-TextValue $ws.Range('B48') 'Cronos'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D48') '0.0516'
Set-TextValue $ws.Range('E48') '  +1.36%  '

Set-TextValue $ws.Range('B49') 'USDD'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue $ws.Range('D49') '1.00'
Set-TextValue $ws.Range('E49') '  +0.07%  '

Set-TextValue $ws.Range('B50') 'EnergySwap'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D50') '7.44'
Set-TextValue $ws.Range('E50') '  +0.54%  '

Set-TextValue $ws.Range('B51') 'Mantle'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D51') '0.396'
Set-TextValue $ws.Range('E51') '  +0.94%  '

